# "Fixed Timer 2 anad updated Unit test protocol sheet"
#
# Reworks the "LCD Driver" test case row (row 6) of the Unit Tests sheet:
#  - Assignee changes from "Mohab/Omar" to just "Omar"
#  - Test steps / expected / actual results are replaced with the new
#    "increment character" LCD test description instead of the old
#    "print A50" description.
# Also moves the sheet's saved selection to C13 and widens column D to fit
# the longer test-steps text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the LCD Driver row (Excel row 6). Cells are written in the same
# left-to-right-ish order the original author typed them in (Expected
# Result, Test Case Steps, Actual Result, then Assignee) so freshly
# introduced shared strings land in that order.
$ws.Range("E6").Value = "LCD should print -> a b c d and so on"
$ws.Range("D6").Value = "Call LCD_sendChar('a') and increment it,Test Upper, lower nibble & 8 bit mode by adjusting config file"
$ws.Range("F6").Value = "LCD printed the expected values"
$ws.Range("C6").Value = "Omar"

# Column D needs to be a little wider now that the test-steps text is longer.
$ws.Columns.Item(4).ColumnWidth = 114.59244791666667

# Move the sheet's active selection to C13.
$null = $ws.Range("C13").Select()
